$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Main title: merge the 5 runs that spell out the old title into a
#    single run with the new wording (Find/Replace across the whole
#    phrase collapses the run boundaries the way the target OOXML
#    expects).
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "DOCUMENTO DE INFORMACIÓN Y AUTORIZACIÓN PARA LA REALIZACIÓN DE LA CIRUGÍA DEL COLESTEATOMA (OÍDO)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "DOCUMENTO DE INFORMACIÓN Y AUTORIZACIÓN PARA LA REALIZACIÓN DE LA CIRUGÍA DE  LA MASTOIDES", 2)

# ------------------------------------------------------------------
# 2. "FECHA DE REALIZACIÓN DEL PROTOCOLO: " -> split the run that held
#    "REALIZACIÓN DEL PROTOCOLO" into two runs: "REALIZACIÓN DEL " and
#    "CONSENTIMIENTO" (same bold/font/size). The surrounding runs
#    ("FECHA DE ", ":" and the trailing space) must stay untouched.
#
#    A plain Find/Replace (or Range.Text=) on this engine re-coalesces
#    every adjacent run sharing the post-edit formatting, which would
#    swallow the neighbouring runs. Briefly nudging Bold off/on around
#    the edited sub-range keeps it from bleeding into its neighbours;
#    reverting it afterwards (a pure formatting no-op, not a text
#    mutation) does not retrigger the coalescing pass.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("REALIZACIÓN DEL PROTOCOLO")
$protoLen = 9                       # Len("PROTOCOLO")
$protoStart = $anchor.End - $protoLen
$protoEnd = $anchor.End

# Give the trailing ":" a momentarily different look than the space
# that follows it, so the two don't get folded back together while the
# "PROTOCOLO" run next to them is being rewritten.
$colon = $d.Range($protoEnd, $protoEnd + 1)
$colon.Bold = 0

# Rewrite PROTOCOLO -> CONSENTIMIENTO, isolated from its neighbours via
# the same Bold nudge.
$proto = $d.Range($protoStart, $protoEnd)
$proto.Bold = 0
$proto.Text = "CONSENTIMIENTO"
$newEnd = $proto.End

# Restore bold on the freshly written "CONSENTIMIENTO" run.
$fixed = $d.Range($protoStart, $newEnd)
$fixed.Bold = 1

# Restore bold on the ":" run (pure formatting — safe, no re-merge).
$colonFixed = $d.Range($newEnd, $newEnd + 1)
$colonFixed.Bold = 1
